$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.638.42"
$ws.Range("E2").Value = "  -5.83%  "
$ws.Range("D3").Value = "1.805.34"
$ws.Range("E3").Value = "  -5.17%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'276.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -9.73%  "
$ws.Range("D7").Value = "'0.5081"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.67%  "
$ws.Range("D8").Value = "'0.3517"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.58%  "
$ws.Range("D9").Value = "'0.06640"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.86%  "
$ws.Range("D10").Value = "'20.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.57%  "
$ws.Range("D11").Value = "'0.8370"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.51%  "
$ws.Range("D12").Value = "'0.07774"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.88%  "
$ws.Range("D13").Value = "1.805.78"
$ws.Range("E13").Value = "  +61.42%  "
$ws.Range("D14").Value = "'5.069"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.21%  "
$ws.Range("D15").Value = "'87.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.66%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "'13.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.36%  "
$ws.Range("D18").Value = "'0.9998"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "'0.000007944"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.34%  "
$ws.Range("D20").Value = "25.694.33"
$ws.Range("E20").Value = "  -5.75%  "
$ws.Range("D21").Value = "'4.711"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.66%  "
$ws.Range("D22").Value = "2.035.04"
$ws.Range("E22").Value = "  +60.69%  "
$ws.Range("D23").Value = "'10.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.43%  "
$ws.Range("D24").Value = "'6.033"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'142.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.97%  "
$ws.Range("D26").Value = "'2.110"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'1.654"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.69%  "
$ws.Range("D28").Value = "'16.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.18%  "
$ws.Range("D29").Value = "'108.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.19%  "
$ws.Range("D30").Value = "'4.314"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -11.12%  "
$ws.Range("D31").Value = "'4.218"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.27%  "
$ws.Range("D32").Value = "'0.08801"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.51%  "
$ws.Range("D33").Value = "'0.04790"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.68%  "
$ws.Range("D34").Value = "'0.7231"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -12.58%  "
$ws.Range("D35").Value = "'1.124"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.72%  "
$ws.Range("D36").Value = "'2.855"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.91%  "
$ws.Range("D37").Value = "'0.9991"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "'3.034"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.57%  "
$ws.Range("D39").Value = "'0.01859"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.02%  "
$ws.Range("D40").Value = "'0.5170"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -12.62%  "
$ws.Range("D41").Value = "'2.295"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -14.41%  "
$ws.Range("D42").Value = "'0.9607"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -11.20%  "
$ws.Range("D43").Value = "'114.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("D44").Value = "'6.181"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.26%  "
$ws.Range("D45").Value = "'8.026"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -13.47%  "
$ws.Range("D46").Value = "'0.9999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "'0.4578"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.38%  "
$ws.Range("D48").Value = "'0.1384"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.45%  "
$ws.Range("D49").Value = "'9.274"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.73%  "
$ws.Range("D50").Value = "'35.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.35%  "
$ws.Range("D51").Value = "'1.491"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.23%  "
